# Add "Area" / "Atotal" columns to the discharge worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"

# Area formulas (column G), rows 2-14
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"
$ws.Range("G4:G11").Formula = "=(D4-D3)*B4/100"
$ws.Range("G12:G14").Formula = "=(D12-D11)*B12/100"

# Total area (column H)
$ws.Range("H2").Formula = "=SUM(G2:G11)"

# Update selection to match the post-edit workbook state
$ws.Range("H2").Select()
